$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Predicted values for the new 't+3' column (column C)
$values = @(
    2,
    -0.3112759201248513,
    0.2068540516619886,
    0.1424338389046362,
    -0.07641264223041008,
    0.1432386406213986,
    -0.4522942319464852,
    -0.2592776833425081,
    -0.4111143807139855,
    0.3583482123923956,
    -0.2395354272017014,
    -0.08679071302812461,
    -0.01165389346983012,
    0.05155204817547648,
    -0.01410836412500775,
    0.3469334748427019,
    0.5975469357316218,
    0.03667612137647214,
    0.3958918984112226,
    0.2938092751582803,
    0.5186086745320461,
    0.3776735756792909,
    -0.06087923480789817,
    4.57384989143239,
    0.4894903269309043,
    0.3649107407635929,
    0.3022116221255643,
    1.066131646617195,
    5.447503492472463,
    0.9548851282734272,
    -0.2129154423199772,
    0.7871732809015993,
    0.9152751613758854,
    -0.6507460697090554,
    0.8138148514943915,
    0.7678300236883154,
    0.741357205273627,
    0.7174454321700235,
    0.581254185786493,
    0.7529239259480895,
    0.5542491414506645,
    0.6901709337605276,
    0.7089146264171076,
    0.6744044547749227,
    0.6583104991077348,
    -1.264221901809358,
    -0.9798357644125127,
    -0.8736975128145139,
    -0.642552574792616,
    -0.04890637255914625,
    -0.8653123345455569,
    -0.8653123345455569,
    -1.097711744551318,
    -0.1870535167737291,
    -0.9995363594736313,
    -0.8846002655561328,
    -0.9340088073217468,
    -1.12794034965257,
    -0.8480402148460283,
    -0.4971111734956855,
    0.3715007413611734,
    -1.208923604921344,
    -0.738311818883726,
    -0.9019909101174564,
    -0.1110487109336713,
    -0.7899089345158966,
    -0.7620704255356806
)

# Loop over the rows and fill column C with the predicted values
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Give the header cell C1 the same formatting as the other header cells (A1, B1)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
